$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "trafo_id" header to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Move the active selection to G6 (cosmetic, matches recorded selection state)
$ws.Range("G6").Select()
